# Natmi following Dr Hou advice
# Recompute the Rtn4-Rtn4rl1 ligand-receptor pair table: every
# (Sending cluster x Target cluster) combination among {ECs, FAPs, M2, sCs}
# now appears (12 rows instead of 8, since Target cluster gains "M2"),
# and every metric column is refreshed with the new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rtn4"
$ws.Range("C2").Value = "Rtn4rl1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 39.66867433333334
$ws.Range("H2").Value = 119.006023
$ws.Range("I2").Value = 0.154574216411057
$ws.Range("J2").Value = 0.1545742164110569
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1028113333333333
$ws.Range("N2").Value = 0.308434
$ws.Range("O2").Value = 0.02963305562291539
$ws.Range("P2").Value = 0.02963305562291539
$ws.Range("Q2").Value = 4.078389299775778
$ws.Range("R2").Value = 36.705503697982
$ws.Range("S2").Value = 0.004580506352777412
$ws.Range("T2").Value = 0.004580506352777411

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rtn4"
$ws.Range("C3").Value = "Rtn4rl1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 39.66867433333334
$ws.Range("H3").Value = 119.006023
$ws.Range("I3").Value = 0.154574216411057
$ws.Range("J3").Value = 0.1545742164110569
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.826632666666667
$ws.Range("N3").Value = 5.479898
$ws.Range("O3").Value = 0.5264858032574322
$ws.Range("P3").Value = 0.5264858032574322
$ws.Range("Q3").Value = 72.46009638062823
$ws.Range("R3").Value = 652.140867425654
$ws.Range("S3").Value = 0.08138113049006349
$ws.Range("T3").Value = 0.08138113049006347

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rtn4"
$ws.Range("C4").Value = "Rtn4rl1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 39.66867433333334
$ws.Range("H4").Value = 119.006023
$ws.Range("I4").Value = 0.154574216411057
$ws.Range("J4").Value = 0.1545742164110569
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.540037333333333
$ws.Range("N4").Value = 4.620112
$ws.Range("O4").Value = 0.4438811411196524
$ws.Range("P4").Value = 0.4438811411196525
$ws.Range("Q4").Value = 61.09123943717511
$ws.Range("R4").Value = 549.8211549345759
$ws.Range("S4").Value = 0.06861257956821606
$ws.Range("T4").Value = 0.06861257956821606

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rtn4"
$ws.Range("C5").Value = "Rtn4rl1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 57.66057933333332
$ws.Range("H5").Value = 172.981738
$ws.Range("I5").Value = 0.2246820449144221
$ws.Range("J5").Value = 0.2246820449144221
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1028113333333333
$ws.Range("N5").Value = 0.308434
$ws.Range("O5").Value = 0.02963305562291539
$ws.Range("P5").Value = 0.02963305562291539
$ws.Range("Q5").Value = 5.928161042032443
$ws.Range("R5").Value = 53.35344937829199
$ws.Range("S5").Value = 0.006658015534419444
$ws.Range("T5").Value = 0.006658015534419444

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rtn4"
$ws.Range("C6").Value = "Rtn4rl1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 57.66057933333332
$ws.Range("H6").Value = 172.981738
$ws.Range("I6").Value = 0.2246820449144221
$ws.Range("J6").Value = 0.2246820449144221
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.826632666666667
$ws.Range("N6").Value = 5.479898
$ws.Range("O6").Value = 0.5264858032574322
$ws.Range("P6").Value = 0.5264858032574322
$ws.Range("Q6").Value = 105.3246977891915
$ws.Range("R6").Value = 947.9222801027239
$ws.Range("S6").Value = 0.118291906894292
$ws.Range("T6").Value = 0.118291906894292

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rtn4"
$ws.Range("C7").Value = "Rtn4rl1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 57.66057933333332
$ws.Range("H7").Value = 172.981738
$ws.Range("I7").Value = 0.2246820449144221
$ws.Range("J7").Value = 0.2246820449144221
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.540037333333333
$ws.Range("N7").Value = 4.620112
$ws.Range("O7").Value = 0.4438811411196524
$ws.Range("P7").Value = 0.4438811411196525
$ws.Range("Q7").Value = 88.79944483496176
$ws.Range("R7").Value = 799.1950035146558
$ws.Range("S7").Value = 0.09973212248571067
$ws.Range("T7").Value = 0.09973212248571067

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Rtn4"
$ws.Range("C8").Value = "Rtn4rl1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 64.993678
$ws.Range("H8").Value = 194.981034
$ws.Range("I8").Value = 0.2532564301015895
$ws.Range("J8").Value = 0.2532564301015895
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1028113333333333
$ws.Range("N8").Value = 0.308434
$ws.Range("O8").Value = 0.02963305562291539
$ws.Range("P8").Value = 0.02963305562291539
$ws.Range("Q8").Value = 6.682086693417333
$ws.Range("R8").Value = 60.13878024075601
$ws.Range("S8").Value = 0.007504761880061385
$ws.Range("T8").Value = 0.007504761880061386

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Rtn4"
$ws.Range("C9").Value = "Rtn4rl1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 64.993678
$ws.Range("H9").Value = 194.981034
$ws.Range("I9").Value = 0.2532564301015895
$ws.Range("J9").Value = 0.2532564301015895
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.826632666666667
$ws.Range("N9").Value = 5.479898
$ws.Range("O9").Value = 0.5264858032574322
$ws.Range("P9").Value = 0.5264858032574322
$ws.Range("Q9").Value = 118.7195753616147
$ws.Range("R9").Value = 1068.476178254532
$ws.Range("S9").Value = 0.1333359150321451
$ws.Range("T9").Value = 0.1333359150321451

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Rtn4"
$ws.Range("C10").Value = "Rtn4rl1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 64.993678
$ws.Range("H10").Value = 194.981034
$ws.Range("I10").Value = 0.2532564301015895
$ws.Range("J10").Value = 0.2532564301015895
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.540037333333333
$ws.Range("N10").Value = 4.620112
$ws.Range("O10").Value = 0.4438811411196524
$ws.Range("P10").Value = 0.4438811411196525
$ws.Range("Q10").Value = 100.0926905506453
$ws.Range("R10").Value = 900.8342149558081
$ws.Range("S10").Value = 0.112415753189383
$ws.Range("T10").Value = 0.1124157531893831

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Rtn4"
$ws.Range("C11").Value = "Rtn4rl1"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 94.308965
$ws.Range("H11").Value = 282.926895
$ws.Range("I11").Value = 0.3674873085729315
$ws.Range("J11").Value = 0.3674873085729314
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1028113333333333
$ws.Range("N11").Value = 0.308434
$ws.Range("O11").Value = 0.02963305562291539
$ws.Range("P11").Value = 0.02963305562291539
$ws.Range("Q11").Value = 9.696030436936665
$ws.Range("R11").Value = 87.26427393243
$ws.Range("S11").Value = 0.01088977185565715
$ws.Range("T11").Value = 0.01088977185565715

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Rtn4"
$ws.Range("C12").Value = "Rtn4rl1"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 94.308965
$ws.Range("H12").Value = 282.926895
$ws.Range("I12").Value = 0.3674873085729315
$ws.Range("J12").Value = 0.3674873085729314
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.826632666666667
$ws.Range("N12").Value = 5.479898
$ws.Range("O12").Value = 0.5264858032574322
$ws.Range("P12").Value = 0.5264858032574322
$ws.Range("Q12").Value = 172.2678362285234
$ws.Range("R12").Value = 1550.41052605671
$ws.Range("S12").Value = 0.1934768508409317
$ws.Range("T12").Value = 0.1934768508409317

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Rtn4"
$ws.Range("C13").Value = "Rtn4rl1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 94.308965
$ws.Range("H13").Value = 282.926895
$ws.Range("I13").Value = 0.3674873085729315
$ws.Range("J13").Value = 0.3674873085729314
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.540037333333333
$ws.Range("N13").Value = 4.620112
$ws.Range("O13").Value = 0.4438811411196524
$ws.Range("P13").Value = 0.4438811411196525
$ws.Range("Q13").Value = 145.2393269680267
$ws.Range("R13").Value = 1307.15394271224
$ws.Range("S13").Value = 0.1631206858763427
$ws.Range("T13").Value = 0.1631206858763427
